# Fix target related bugs
# NoSales.xlsx: replace the stale "no-sales" product listing (Esoral/Losectil/
# Rabifast/Softi) with the corrected target-brand listing (Aldorin/Cardoneb/
# Cardovan/Dialon/GLIKAZID/Irbes/Ligazid/Lipicon/Pivasta/Sitazid) and drop the
# trailing row that no longer belongs in the list (21 data rows -> 20).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colA = @(4, 17, 17, 18, 18, 18, 18, 26, 52, 57, 68, 68, 68, 70, 70, 70, 70, 104, 123, 123)
$colB = @("Aldorin", "Cardoneb", "Cardoneb", "Cardovan", "Cardovan", "Cardovan", "Cardovan", "Dialon", "GLIKAZID", "Irbes", "Ligazid", "Ligazid", "Ligazid", "Lipicon", "Lipicon", "Lipicon", "Lipicon", "Pivasta", "Sitazid", "Sitazid")
$colC = @(1, 2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20)
$colD = @("Aldorin 50mg Tablet - 24's", "Cardoneb 5 FC Tablet 30's", "Cardoneb 2.5 FC Tablet 30's", "Cardovan Plus 80/12.5 Tablet 30's", "Cardovan 160mg Tablet 30's", "Cardovan 80mg Tablet 30's", "Cardovan Plus 160/12.5 Tablet 30's", "Dialon 4mg Tablet", "Glikazid 80mg Tablet 30's", "Irbes 75mg Tablet", "Ligazid 5mg Tablet 20's", "Ligazid 5mg Tablet 10's", "Ligazid M 2.5/500", "Lipicon 20mg Tablet - 20's", "Lipicon 10mg Tablet Container 30's", "Lipicon 40mg Tablet - 10's", "Lipicon 10mg Tablet - 40's", "Pivasta 2mg Tablet 20's", "Sitazid 50mg Tablet 20's", "Sitazid 100mg Tablet 10's")
$colE = @("24's", "30's", "30's", "30's", "30's", "30's", "30's", "20'S", "30's", "50 's", "20's", "10's", "20's", "20 's", "30's", "10 's", "40 's", "20's", "20's", "10's")

# Drop the last data row (was row 22, BSL 165 / Softi) -- the target sheet
# only has 20 data rows (rows 2-21) instead of 21.
$ws.Rows.Item(22).Delete()

$rowCount = $colA.Length

# Blank out the text columns first so every stale shared string (old BRAND /
# Item Name / UOM values) is fully dereferenced and dropped from the table
# before the new values are written back in.
for ($r = 2; $r -le ($rowCount + 1); $r++) {
    $ws.Cells.Item($r, 2).Value = ""
    $ws.Cells.Item($r, 4).Value = ""
    $ws.Cells.Item($r, 5).Value = ""
}

# Now write column-by-column so the shared-string table is rebuilt in the
# same grouped order the source file uses (all BRAND values, then all Item
# Name values, then all UOM values).
for ($i = 0; $i -lt $rowCount; $i++) {
    $ws.Cells.Item(2 + $i, 1).Value = $colA[$i]
}
for ($i = 0; $i -lt $rowCount; $i++) {
    $ws.Cells.Item(2 + $i, 2).Value = $colB[$i]
}
for ($i = 0; $i -lt $rowCount; $i++) {
    $ws.Cells.Item(2 + $i, 3).Value = $colC[$i]
}
for ($i = 0; $i -lt $rowCount; $i++) {
    $ws.Cells.Item(2 + $i, 4).Value = $colD[$i]
}
for ($i = 0; $i -lt $rowCount; $i++) {
    $ws.Cells.Item(2 + $i, 5).Value = $colE[$i]
}

